$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Pass 1: write new shared-string text cells in exact allocation order ---
$ws.Range("A261").Value = "b9861_d180530_TT2_SS01"
$ws.Range("A262").Value = "b9861_d180530_TT2_SS02"
$ws.Range("A263").Value = "b9861_d180530_TT2_SS03"
$ws.Range("A264").Value = "b9861_d180531_TT3_SS01"
$ws.Range("A265").Value = "b9861_d180531_TT4_SS01"
$ws.Range("A266").Value = "b9861_d180531_TT4_SS02"
$ws.Range("A267").Value = "b9861_d180601_TT1_SS01"
$ws.Range("A268").Value = "b9861_d180601_TT1_SS02"
$ws.Range("A269").Value = "b9861_d180601_TT1_SS03"
$ws.Range("A271").Value = "b9861_d180601_TT3_SS02"
$ws.Range("A270").Value = "b9861_d180601_TT3_SS01"
$ws.Range("A272").Value = "b9861_d180601_TT3_SS03"
$ws.Range("A273").Value = "b9861_d180601_TT3_SS04"
$ws.Range("A274").Value = "b9861_d180601_TT3_SS05"
$ws.Range("A275").Value = "b9861_d180603_TT1_SS01"
$ws.Range("G275").Value = "lots of flight artifacts in that TT"
$ws.Range("A276").Value = "b9861_d180603_TT2_SS01"
$ws.Range("A277").Value = "b9861_d180603_TT3_SS01"
$ws.Range("A278").Value = "b9861_d180603_TT3_SS02"
$ws.Range("A279").Value = "b9861_d180604_TT3_SS01"
$ws.Range("A280").Value = "b9861_d180604_TT3_SS02"
$ws.Range("A281").Value = "b9861_d180604_TT3_SS03"
$ws.Range("A282").Value = "b9861_d180604_TT3_SS04"
$ws.Range("A283").Value = "b9861_d180604_TT3_SS05"
$ws.Range("A284").Value = "b9861_d180605_TT3_SS01"
$ws.Range("A285").Value = "b9861_d180605_TT3_SS02"
$ws.Range("A287").Value = "b9861_d180606_TT3_SS02"
$ws.Range("A286").Value = "b9861_d180606_TT3_SS01"
$ws.Range("A288").Value = "b9861_d180606_TT3_SS03"
$ws.Range("A289").Value = "b9861_d180606_TT3_SS04"
$ws.Range("A290").Value = "b9861_d180606_TT3_SS05"
$ws.Range("A291").Value = "b9861_d180606_TT3_SS06"
$ws.Range("A292").Value = "b9861_d180607_TT3_SS01"
$ws.Range("A293").Value = "b9861_d180609_TT3_SS01"
$ws.Range("A294").Value = "b9861_d180609_TT3_SS02"
$ws.Range("A295").Value = "b9861_d180610_TT3_SS01"

# --- Pass 2: write G-column cells that reuse EXISTING shared strings ---
$ws.Range("G262").Value = "not well isolated"
$ws.Range("G265").Value = "partially stable"
$ws.Range("G266").Value = "multi-unit close to the noise; partially stable"
$ws.Range("G268").Value = "not well isolated"
$ws.Range("G269").Value = "not well isolated"
$ws.Range("G271").Value = "multi-unit close to the noise"
$ws.Range("G272").Value = "multi-unit close to the noise"
$ws.Range("G273").Value = "multi-unit close to the noise"
$ws.Range("G274").Value = "multi-unit close to the noise"
$ws.Range("G276").Value = "multi-unit close to the noise"
$ws.Range("G278").Value = "multi-unit close to the noise"
$ws.Range("G279").Value = "multi-unit close to the noise"
$ws.Range("G280").Value = "multi-unit close to the noise"
$ws.Range("G281").Value = "multi-unit close to the noise"
$ws.Range("G282").Value = "multi-unit close to the noise"
$ws.Range("G283").Value = "multi-unit close to the noise"
$ws.Range("G284").Value = "multi-unit close to the noise"
$ws.Range("G285").Value = "multi-unit close to the noise"
$ws.Range("G289").Value = "multi-unit close to the noise"
$ws.Range("G290").Value = "multi-unit close to the noise"
$ws.Range("G291").Value = "multi-unit close to the noise"
$ws.Range("G292").Value = "multi-unit close to the noise"
$ws.Range("G293").Value = "multi-unit close to the noise"
$ws.Range("G294").Value = "multi-unit close to the noise"
$ws.Range("G295").Value = "multi-unit close to the noise"

# --- Pass 3: numeric columns B-F for every new row ---
$ws.Range("B261").Value = 9861
$ws.Range("C261").Value = 43250
$ws.Range("D261").Value = 2
$ws.Range("E261").Value = 1
$ws.Range("F261").Value = 2
$ws.Range("B262").Value = 9861
$ws.Range("C262").Value = 43250
$ws.Range("D262").Value = 2
$ws.Range("E262").Value = 2
$ws.Range("F262").Value = 1
$ws.Range("B263").Value = 9861
$ws.Range("C263").Value = 43250
$ws.Range("D263").Value = 2
$ws.Range("E263").Value = 3
$ws.Range("F263").Value = 2
$ws.Range("B264").Value = 9861
$ws.Range("C264").Value = 43251
$ws.Range("D264").Value = 3
$ws.Range("E264").Value = 1
$ws.Range("F264").Value = 2
$ws.Range("B265").Value = 9861
$ws.Range("C265").Value = 43251
$ws.Range("D265").Value = 4
$ws.Range("E265").Value = 1
$ws.Range("F265").Value = 2
$ws.Range("B266").Value = 9861
$ws.Range("C266").Value = 43251
$ws.Range("D266").Value = 4
$ws.Range("E266").Value = 2
$ws.Range("F266").Value = 1
$ws.Range("B267").Value = 9861
$ws.Range("C267").Value = 43252
$ws.Range("D267").Value = 1
$ws.Range("E267").Value = 1
$ws.Range("F267").Value = 2
$ws.Range("B268").Value = 9861
$ws.Range("C268").Value = 43252
$ws.Range("D268").Value = 1
$ws.Range("E268").Value = 2
$ws.Range("F268").Value = 1
$ws.Range("B269").Value = 9861
$ws.Range("C269").Value = 43252
$ws.Range("D269").Value = 1
$ws.Range("E269").Value = 3
$ws.Range("F269").Value = 1
$ws.Range("B270").Value = 9861
$ws.Range("C270").Value = 43252
$ws.Range("D270").Value = 3
$ws.Range("E270").Value = 1
$ws.Range("F270").Value = 2
$ws.Range("B271").Value = 9861
$ws.Range("C271").Value = 43252
$ws.Range("D271").Value = 3
$ws.Range("E271").Value = 2
$ws.Range("F271").Value = 1
$ws.Range("B272").Value = 9861
$ws.Range("C272").Value = 43252
$ws.Range("D272").Value = 3
$ws.Range("E272").Value = 3
$ws.Range("F272").Value = 1
$ws.Range("B273").Value = 9861
$ws.Range("C273").Value = 43252
$ws.Range("D273").Value = 3
$ws.Range("E273").Value = 4
$ws.Range("F273").Value = 1
$ws.Range("B274").Value = 9861
$ws.Range("C274").Value = 43252
$ws.Range("D274").Value = 3
$ws.Range("E274").Value = 5
$ws.Range("F274").Value = 1
$ws.Range("B275").Value = 9861
$ws.Range("C275").Value = 43254
$ws.Range("D275").Value = 1
$ws.Range("E275").Value = 1
$ws.Range("F275").Value = 2
$ws.Range("B276").Value = 9861
$ws.Range("C276").Value = 43254
$ws.Range("D276").Value = 2
$ws.Range("E276").Value = 1
$ws.Range("F276").Value = 1
$ws.Range("B277").Value = 9861
$ws.Range("C277").Value = 43254
$ws.Range("D277").Value = 3
$ws.Range("E277").Value = 1
$ws.Range("F277").Value = 2
$ws.Range("B278").Value = 9861
$ws.Range("C278").Value = 43254
$ws.Range("D278").Value = 3
$ws.Range("E278").Value = 2
$ws.Range("F278").Value = 1
$ws.Range("B279").Value = 9861
$ws.Range("C279").Value = 43255
$ws.Range("D279").Value = 3
$ws.Range("E279").Value = 1
$ws.Range("F279").Value = 1
$ws.Range("B280").Value = 9861
$ws.Range("C280").Value = 43255
$ws.Range("D280").Value = 3
$ws.Range("E280").Value = 2
$ws.Range("F280").Value = 1
$ws.Range("B281").Value = 9861
$ws.Range("C281").Value = 43255
$ws.Range("D281").Value = 3
$ws.Range("E281").Value = 3
$ws.Range("F281").Value = 1
$ws.Range("B282").Value = 9861
$ws.Range("C282").Value = 43255
$ws.Range("D282").Value = 3
$ws.Range("E282").Value = 4
$ws.Range("F282").Value = 1
$ws.Range("B283").Value = 9861
$ws.Range("C283").Value = 43255
$ws.Range("D283").Value = 3
$ws.Range("E283").Value = 5
$ws.Range("F283").Value = 1
$ws.Range("B284").Value = 9861
$ws.Range("C284").Value = 43256
$ws.Range("D284").Value = 3
$ws.Range("E284").Value = 1
$ws.Range("F284").Value = 1
$ws.Range("B285").Value = 9861
$ws.Range("C285").Value = 43256
$ws.Range("D285").Value = 3
$ws.Range("E285").Value = 2
$ws.Range("F285").Value = 1
$ws.Range("B286").Value = 9861
$ws.Range("C286").Value = 43257
$ws.Range("D286").Value = 3
$ws.Range("E286").Value = 1
$ws.Range("F286").Value = 2
$ws.Range("B287").Value = 9861
$ws.Range("C287").Value = 43257
$ws.Range("D287").Value = 3
$ws.Range("E287").Value = 2
$ws.Range("F287").Value = 2
$ws.Range("B288").Value = 9861
$ws.Range("C288").Value = 43257
$ws.Range("D288").Value = 3
$ws.Range("E288").Value = 3
$ws.Range("F288").Value = 2
$ws.Range("B289").Value = 9861
$ws.Range("C289").Value = 43257
$ws.Range("D289").Value = 3
$ws.Range("E289").Value = 4
$ws.Range("F289").Value = 1
$ws.Range("B290").Value = 9861
$ws.Range("C290").Value = 43257
$ws.Range("D290").Value = 3
$ws.Range("E290").Value = 5
$ws.Range("F290").Value = 1
$ws.Range("B291").Value = 9861
$ws.Range("C291").Value = 43257
$ws.Range("D291").Value = 3
$ws.Range("E291").Value = 6
$ws.Range("F291").Value = 1
$ws.Range("B292").Value = 9861
$ws.Range("C292").Value = 43258
$ws.Range("D292").Value = 3
$ws.Range("E292").Value = 1
$ws.Range("F292").Value = 1
$ws.Range("B293").Value = 9861
$ws.Range("C293").Value = 43260
$ws.Range("D293").Value = 3
$ws.Range("E293").Value = 1
$ws.Range("F293").Value = 1
$ws.Range("B294").Value = 9861
$ws.Range("C294").Value = 43260
$ws.Range("D294").Value = 3
$ws.Range("E294").Value = 2
$ws.Range("F294").Value = 1
$ws.Range("B295").Value = 9861
$ws.Range("C295").Value = 43261
$ws.Range("D295").Value = 3
$ws.Range("E295").Value = 1
$ws.Range("F295").Value = 1

# --- Pass 4: number format for column C (date) on new rows ---
$ws.Range("C261:C295").NumberFormat = "m/d/yy"

# --- Pass 5: highlight fill for G265 and G266 (matches existing yellow-highlight style) ---
$ws.Range("G265:G266").Interior.Color = 65535
